$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "male"
$ws.Range("C3").Value = "singh"
$ws.Range("B3").Value = "ravinder"
$ws.Range("A1").Value = "User register details "
$ws.Range("A2").Value = "Gender:"
$ws.Range("B2").Value = "First name:"
$ws.Range("C2").Value = "Last name:"
$ws.Range("D2").Value = "Password:"
$ws.Range("D3").Value = "qwerty123"
$ws.Range("A4").Value = "Product to search"
$ws.Range("A5").Value = "Product name:"
$ws.Range("A6").Value = "Apple MacBook"
$ws.Range("B6").Value = "Apple MacBook Pro 13-inch"
$ws.Range("B5").Value = "Product searched verify:"

$ws.Range("A1").Interior.Color = 5296274
$ws.Range("A1").Font.Bold = $true

$ws.Range("A1").Copy()
$ws.Range("A4").PasteSpecial(-4122)
